$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 959
$ws.Range("B2").Value = 904
$ws.Range("C2").Value = 904
$ws.Range("D2").Value = 904
$ws.Range("E2").Value = 939
$ws.Range("F2").Value = 965
$ws.Range("G2").Value = 955
$ws.Range("H2").Value = 974
